# fix(publipostage): Correct status name
#
# Renames the "bleu" status label to "noir", and rewords the
# "statut_name" descriptions to use "postés ... publiés" phrasing
# instead of "et / ou publication posté", across every row of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "bleu" = "noir"
    "pas de résultat ni de publication" = "pas de résultat postés ni publiés"
    "résultat et / ou publication posté" = "résultat postés ou publiés"
    "résultat et / ou publication posté dans les 36 mois" = "résultat postés ou publiés dans les 36 mois"
    "résultat et / ou publication posté dans les 12 mois" = "résultat postés ou publiés dans les 12 mois"
}

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $current = $cell.Value()
        if ($current -ne $null -and $replacements.ContainsKey($current)) {
            $cell.Value = $replacements[$current]
        }
    }
}
